$d = $word.ActiveDocument

# 1. Replace {{#image data}} with {{#image generationChart}}
$r1 = $d.Content
$r1.Find.Execute("{{#image data}}", $true, $false, $false, $false, $false,
                  $true, 1, $false, "{{#image generationChart}}", 2)

Write-Host "Done"
